$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "someone"
$ws.Range("B2").Value = "some code"
$ws.Range("C2").Value = "some title"
$ws.Range("D2").Value = "some type"
$ws.Range("E2").Value = "some hrs"
$ws.Range("F2").Value = "some tut hours"
$ws.Range("G2").Value = "some pracs"
$ws.Range("H2").Value = "some projs"
$ws.Range("I2").Value = "some credits"
$ws.Range("J2").Value = "some id"
$ws.Range("K2").Value = "some room"
$ws.Range("L2").Value = "some slot"
$ws.Range("M2").Value = "some erpid"
$ws.Range("N2").Value = "some emp"
$ws.Range("O2").Value = "some school"
$ws.Range("P2").Value = "some mode"
$ws.Range("Q2").Value = "some"
$ws.Range("R2").Value = "some"
$ws.Range("S2").Value = "some"
$ws.Range("T2").Value = "some"
$ws.Range("U2").Value = "some"
$ws.Range("V2").Value = "some"

$ws.Range("V2").Select()
